$d = $word.ActiveDocument

function Replace-Exact($findText) {
    $d.Content.Find.Execute(
        $findText, $true, $false, $false, $false, $false,
        $true, 1, $false, $findText, 2
    ) | Out-Null
}

# --- Collapse runs that were only split apart by now-removed w:proofErr
#     (spelling/grammar check) markers back into single runs, by doing a
#     same-text Find & Replace across each split span. ---

Replace-Exact("Having at least one of your own hostage or converted law enforcer makes you regenerate ")

Replace-Exact(" armor for each converted enemy up to ")

Replace-Exact("You can now take special enemies hostage and convert them.")

Replace-Exact(" more armor")

Replace-Exact("armor for each successful headshot. Cooldown is reduced ")

Replace-Exact("The movement speed penalty of armor is reduced by ")

Replace-Exact("You have 30% chance to enter 'Bulletstorm' for ")

Replace-Exact("Increases your armor recovery rate by ")

Replace-Exact("You can now bring the normal amount of secondary deployables with you")

Replace-Exact(" faster. You can now hack keycard security panels")

Replace-Exact("The effect persists for 9 seconds after your armor has recovered.")

Replace-Exact("When your armor breaks you gain a ")

Replace-Exact("You now have 9001% more bleedout health")

# --- Actual content edit: insert "mark" into the Hitman AA description ---

$d.Content.Find.Execute(
    "Outside of stealth automatically units who near you.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Outside of stealth automatically mark units who near you.",
    2
) | Out-Null
